$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @'
questions = [
    {
        "title": "As a social media marketing manager, you are planning a campaign targeting 18- to 35-year-old females in the US. Your goal is maximum exposure to the largest audience.Which social media channel would be the most effective for your campaign?",
        "ques_type": 2,
        "options": [
            "Pinterest",
            "Instagram",
            "Facebook",
            "Twitch",
            "Twitter",
            "LinkedIn"
        ],
        "score": "Instagram"
    },
    {
        "title": "As an Amazon PPC manager, you are overseeing an ongoing campaign and need to analyze its performance for optimization.What is the appropriate frequency to review your results?",
        "ques_type": 15,
        "options": [
            "Every four hours in real time",
            "At least once a day in real time",
            "At least twice a day in real time",
            "Every morning, looking at the previous day",
            "Every seven days"
        ],
        "score": [
            "Every morning, looking at the previous day",
            "Every seven days"
        ]
    },
    {
        "title": "You are a marketing strategist at a tech firm that specializes in selling complex software solutions to other businesses. Your manager has asked you to recommend the best primary paid advertising channel for the company.Which platform would you suggest as the most effective for this purpose?",
        "ques_type": 2,
        "options": [
            "Instagram",
            "Google Ads",
            "LinkedIn",
            "YouTube",
            "Facebook Ads"
        ],
        "score": "LinkedIn"
    },
    {
        "title": "You are a digital marketing specialist managing search engine ads for your company. Your supervisor is questioning the relationship between quality score and the cost of search engine ads, and they want your input.How should you respond?",
        "ques_type": 2,
        "options": [
            "A higher quality score increases the ad prices because of increased competition.",
            "A higher quality score doesn't affect ad prices it only improves ad ranking.",
            "A higher quality score reduces ad prices due to improved relevance and user experience.",
            "A higher quality score causes unpredictable fluctuations in ad prices.",
            "A higher quality score increases ad prices because it leads to more impressions."
        ],
        "score": "A higher quality score reduces ad prices due to improved relevance and user experience."
    }
]
'@

# Remove the old header row (A1 = 0 with bold/bordered style); this shifts
# the question-data row (old A2) up into A1 and drops the unused style.
$ws.Rows.Item(1).Delete()

# Replace the shared-string content with the pretty-printed JSON form.
$ws.Range("A1").Value2 = $newText

# Setting a multi-line value auto-expands the row height; restore the
# sheet's natural (non-custom) row height so the row tag stays plain.
$ws.Rows.Item(1).AutoFit()
